$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 117.57143
$ws.Range("I33").Value = 144.33333
$ws.Range("J33").Value = 69.40000000000001
$ws.Range("K33").Value = 144.33333
$ws.Range("L33").Value = 69.40000000000001
$ws.Range("M33").Value = 84.66667000000001
$ws.Range("N33").Value = -527.4
$ws.Range("H38").Value = 12474.875
$ws.Range("J38").Value = 15600
$ws.Range("L38").Value = 46800
$ws.Range("N38").Value = -47544
$ws.Range("H40").Value = 27794226
$ws.Range("J40").Value = 38474652
$ws.Range("L40").Value = 38474652
$ws.Range("N40").Value = -38475002
$ws.Range("H92").Value = 38461904
$ws.Range("I92").Value = 45454890
$ws.Range("K92").Value = 45454890
$ws.Range("M92").Value = -45453642
$ws.Range("H98").Value = 1681
$ws.Range("I98").Value = 1561.4445
$ws.Range("J98").Value = 2111.4
$ws.Range("K98").Value = 1561.4445
$ws.Range("L98").Value = 2111.4
$ws.Range("M98").Value = -63.44450000000006
$ws.Range("N98").Value = -5107.4
$ws.Range("H100").Value = 1600
$ws.Range("I100").Value = 1000
$ws.Range("K100").Value = 1000
$ws.Range("M100").Value = -459
$ws.Range("H101").Value = 23966.8
$ws.Range("I101").Value = 3574.5
$ws.Range("K101").Value = 10723.5
$ws.Range("M101").Value = -9101.5
$ws.Range("H112").Value = 4568.483
$ws.Range("I112").Value = 3698
$ws.Range("J112").Value = 4599.5713
$ws.Range("K112").Value = 11094
$ws.Range("L112").Value = 13798.7139
$ws.Range("M112").Value = -9986
$ws.Range("N112").Value = -16014.7139
$ws.Range("H122").Value = 1681
$ws.Range("I122").Value = 1561.4445
$ws.Range("J122").Value = 2111.4
$ws.Range("K122").Value = 4684.333500000001
$ws.Range("L122").Value = 6334.200000000001
$ws.Range("M122").Value = -2234.333500000001
$ws.Range("N122").Value = -11234.2
$ws.Range("H132").Value = 110794.09
$ws.Range("I132").Value = 302643.5
$ws.Range("K132").Value = 907930.5
$ws.Range("M132").Value = -905400.5
$ws.Range("H141").Value = 4207.185
$ws.Range("I141").Value = 4151.875
$ws.Range("J141").Value = 4649.6665
$ws.Range("K141").Value = 12455.625
$ws.Range("L141").Value = 13948.9995
$ws.Range("M141").Value = -7275.625
$ws.Range("N141").Value = -24308.9995

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 651.8333
$ws.Range("I5").Value = 764.2
$ws.Range("K5").Value = 764.2
$ws.Range("M5").Value = -652.2
$ws.Range("H74").Value = 13895244
$ws.Range("I74").Value = 125000500
$ws.Range("K74").Value = 125000500
$ws.Range("M74").Value = -124999626
$ws.Range("H77").Value = 13895244
$ws.Range("I77").Value = 125000500
$ws.Range("K77").Value = 625002500
$ws.Range("M77").Value = -624998132
$ws.Range("H102").Value = 416728.3
$ws.Range("J102").Value = 1861.8
$ws.Range("L102").Value = 1861.8
$ws.Range("N102").Value = -5105.8
$ws.Range("H132").Value = 10750.25
$ws.Range("I132").Value = 11217.807
$ws.Range("J132").Value = 1866.6666
$ws.Range("K132").Value = 33653.421
$ws.Range("L132").Value = 5599.9998
$ws.Range("M132").Value = -31123.421
$ws.Range("N132").Value = -10659.9998
$ws.Range("H135").Value = 90214.5
$ws.Range("J135").Value = 90214.5
$ws.Range("L135").Value = 90214.5
$ws.Range("N135").Value = -100354.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 651.8333
$ws.Range("I4").Value = 764.2
$ws.Range("K4").Value = 764.2
$ws.Range("M4").Value = -649.2
$ws.Range("H94").Value = 721568
$ws.Range("I94").Value = 806334.9
$ws.Range("J94").Value = 1049.5
$ws.Range("K94").Value = 806334.9
$ws.Range("L94").Value = 1049.5
$ws.Range("M94").Value = -805883.9
$ws.Range("N94").Value = -1951.5
$ws.Range("H99").Value = 102084600
$ws.Range("I99").Value = 255209840
$ws.Range("J99").Value = 1115.3334
$ws.Range("K99").Value = 255209840
$ws.Range("L99").Value = 1115.3334
$ws.Range("M99").Value = -255208342
$ws.Range("N99").Value = -4111.3334
$ws.Range("H105").Value = 150001890
$ws.Range("I105").Value = 214287410
$ws.Range("J105").Value = 2333.3333
$ws.Range("K105").Value = 214287410
$ws.Range("L105").Value = 2333.3333
$ws.Range("M105").Value = -214285663
$ws.Range("N105").Value = -5827.3333
$ws.Range("H134").Value = 4390.784
$ws.Range("I134").Value = 1517.64
$ws.Range("J134").Value = 7153.423
$ws.Range("K134").Value = 4552.92
$ws.Range("L134").Value = 21460.269
$ws.Range("M134").Value = -2017.92
$ws.Range("N134").Value = -26530.269

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1108
$ws.Range("J5").Value = 1108
$ws.Range("L5").Value = 1108
$ws.Range("N5").Value = -1332
$ws.Range("H58").Value = 1819413.9
$ws.Range("I58").Value = 2223228
$ws.Range("J58").Value = 2249.5
$ws.Range("K58").Value = 2223228
$ws.Range("L58").Value = 2249.5
$ws.Range("M58").Value = -2223025
$ws.Range("N58").Value = -2655.5
$ws.Range("H132").Value = 49392028
$ws.Range("I132").Value = 60608110
$ws.Range("J132").Value = 41265
$ws.Range("K132").Value = 181824330
$ws.Range("L132").Value = 123795
$ws.Range("M132").Value = -181821800
$ws.Range("N132").Value = -128855
$ws.Range("H134").Value = 1471.8096
$ws.Range("I134").Value = 1430.75
$ws.Range("J134").Value = 1603.2
$ws.Range("K134").Value = 4292.25
$ws.Range("L134").Value = 4809.6
$ws.Range("M134").Value = -1757.25
$ws.Range("N134").Value = -9879.6
$ws.Range("H136").Value = 1819413.9
$ws.Range("I136").Value = 2223228
$ws.Range("J136").Value = 2249.5
$ws.Range("K136").Value = 6669684
$ws.Range("L136").Value = 6748.5
$ws.Range("M136").Value = -6667134
$ws.Range("N136").Value = -11848.5
$ws.Range("H138").Value = 83796
$ws.Range("I138").Value = 66332.336
$ws.Range("J138").Value = 101259.664
$ws.Range("K138").Value = 66332.336
$ws.Range("L138").Value = 101259.664
$ws.Range("M138").Value = -61192.336
$ws.Range("N138").Value = -111539.664

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 704.44446
$ws.Range("I5").Value = 431.05884
$ws.Range("K5").Value = 1293.17652
$ws.Range("M5").Value = -1181.17652
$ws.Range("H63").Value = 24147
$ws.Range("I63").Value = 23995.5
$ws.Range("J63").Value = 24222.75
$ws.Range("K63").Value = 71986.5
$ws.Range("L63").Value = 72668.25
$ws.Range("M63").Value = -71237.5
$ws.Range("N63").Value = -74166.25
$ws.Range("H66").Value = 24147
$ws.Range("I66").Value = 23995.5
$ws.Range("J66").Value = 24222.75
$ws.Range("K66").Value = 215959.5
$ws.Range("L66").Value = 218004.75
$ws.Range("M66").Value = -212215.5
$ws.Range("N66").Value = -225492.75
$ws.Range("H131").Value = 15156310
$ws.Range("J131").Value = 6433092
$ws.Range("L131").Value = 19299276
$ws.Range("N131").Value = -19309356
$ws.Range("H135").Value = 704.44446
$ws.Range("I135").Value = 431.05884
$ws.Range("K135").Value = 3879.52956
$ws.Range("M135").Value = -1344.52956

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 32001060
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("N3").Value = -5232
$ws.Range("H80").Value = 1286628
$ws.Range("I80").Value = 2380493.5
$ws.Range("J80").Value = 10451.667
$ws.Range("K80").Value = 2380493.5
$ws.Range("L80").Value = 10451.667
$ws.Range("M80").Value = -2379495.5
$ws.Range("N80").Value = -12447.667
$ws.Range("H83").Value = 1286628
$ws.Range("I83").Value = 2380493.5
$ws.Range("J83").Value = 10451.667
$ws.Range("K83").Value = 11902467.5
$ws.Range("L83").Value = 52258.335
$ws.Range("M83").Value = -11897475.5
$ws.Range("N83").Value = -62242.335
$ws.Range("H113").Value = 33350100
$ws.Range("J113").Value = 9997.5
$ws.Range("L113").Value = 9997.5
$ws.Range("N113").Value = -14337.5
$ws.Range("H122").Value = 292943.75
$ws.Range("I122").Value = 395360.44
$ws.Range("J122").Value = 6177
$ws.Range("K122").Value = 1186081.32
$ws.Range("L122").Value = 18531
$ws.Range("M122").Value = -1183631.32
$ws.Range("N122").Value = -23431
$ws.Range("H132").Value = 2499.6445
$ws.Range("I132").Value = 2292.7441
$ws.Range("J132").Value = 6948
$ws.Range("K132").Value = 6878.2323
$ws.Range("L132").Value = 20844
$ws.Range("M132").Value = -4348.2323
$ws.Range("N132").Value = -25904

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 4871.6665
$ws.Range("I31").Value = 15
$ws.Range("J31").Value = 7300
$ws.Range("K31").Value = 15
$ws.Range("L31").Value = 7300
$ws.Range("M31").Value = 233
$ws.Range("N31").Value = -7796
$ws.Range("H46").Value = 5010.684
$ws.Range("J46").Value = 5785.2812
$ws.Range("L46").Value = 5785.2812
$ws.Range("N46").Value = -6161.2812
$ws.Range("H55").Value = 314.9
$ws.Range("I55").Value = 341.8125
$ws.Range("J55").Value = 207.25
$ws.Range("K55").Value = 341.8125
$ws.Range("L55").Value = 207.25
$ws.Range("M55").Value = -168.8125
$ws.Range("N55").Value = -553.25
$ws.Range("H88").Value = 37687.25
$ws.Range("I88").Value = 16916.334
$ws.Range("J88").Value = 100000
$ws.Range("K88").Value = 16916.334
$ws.Range("L88").Value = 100000
$ws.Range("M88").Value = -16488.334
$ws.Range("N88").Value = -100856
$ws.Range("H91").Value = 37687.25
$ws.Range("I91").Value = 16916.334
$ws.Range("J91").Value = 100000
$ws.Range("K91").Value = 16916.334
$ws.Range("L91").Value = 100000
$ws.Range("M91").Value = -15434.334
$ws.Range("N91").Value = -102964
$ws.Range("H93").Value = 2213.963
$ws.Range("I93").Value = 2213.7
$ws.Range("K93").Value = 2213.7
$ws.Range("M93").Value = -965.6999999999998
$ws.Range("H132").Value = 3995
$ws.Range("I132").Value = 3995
$ws.Range("K132").Value = 11985
$ws.Range("M132").Value = -9455
$ws.Range("H136").Value = 3992.08

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 16672867
$ws.Range("I132").Value = 7429
$ws.Range("K132").Value = 22287
$ws.Range("M132").Value = -19757
